# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gets three new trailing columns (H, I, J):
#   date / legislator_name / legislator_id
# populated for every existing data row with the same values:
#   2013-11-08 / 葉宜津 / 855

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Copy the existing header cell's formatting (G1, bold + border) onto the
# new header cells H1:J1, then fill in the new header labels.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Copy the existing data cell's formatting (G2) onto the new data cells
# H2:J4 so the added columns match the look of the rest of the table.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2:J4").PasteSpecial(-4122) | Out-Null

# The date must be stored as plain text "2013-11-08", not an Excel date
# serial number, so force a text format before assigning the value.
$ws.Range("H2:H4").NumberFormat = "@"
$ws.Range("H2").Value = "2013-11-08"
$ws.Range("H3").Value = "2013-11-08"
$ws.Range("H4").Value = "2013-11-08"

$ws.Range("I2").Value = "葉宜津"
$ws.Range("I3").Value = "葉宜津"
$ws.Range("I4").Value = "葉宜津"

$ws.Range("J2").Value = 855
$ws.Range("J3").Value = 855
$ws.Range("J4").Value = 855

$excel.CutCopyMode = 0
